# Update the header date
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-08-05 Tuesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-08-06 Wednesday", 2)

# Update the multiplication problems in the table, row by row / cell by cell
# so duplicate values (e.g. "405x4=") are not ambiguous.
$tbl = $d.Tables.Item(1)

$rowMap = @{
    1  = @("335×8=", "810×2=", "664×2=", "596×5=", "369×2=")
    5  = @("824×3=", "296×8=", "153×9=", "109×7=", "435×8=")
    10 = @("946×3=", "899×7=", "612×2=", "322×8=", "810×8=")
    15 = @("965×5=", "250×7=", "627×9=", "232×9=", "225×8=")
    20 = @("279×4=", "704×6=", "267×7=", "866×7=", "742×7=")
}

foreach ($rowIndex in $rowMap.Keys) {
    $values = $rowMap[$rowIndex]
    $row = $tbl.Rows.Item($rowIndex)
    for ($c = 1; $c -le 5; $c++) {
        $cell = $row.Cells.Item($c)
        $cellRange = $cell.Range
        $cellRange.MoveEnd(1, -1)
        $cellRange.Text = $values[$c - 1]
    }
}
